# Auto-generated Excel COM-interop edit script
# Updates the cryptos price/volume table to the latest scraped values
# (GitHub Actions scheduled data refresh).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text would otherwise be auto-parsed as a number by Excel
# (e.g. '12.05' or '1.00') are written through a text-formatted round trip
# so they stay plain text like the rest of the sheet, then the temporary
# '@' number format is cleared again so no visible style change remains.
$textForced = New-Object System.Collections.Generic.List[object]

# Row 2
$ws.Range("D2").Value = "66.495.54"
$ws.Range("E2").Value = "  -3.57%  "

# Row 3
$ws.Range("D3").Value = "3.555.37"
$ws.Range("E3").Value = "  -4.46%  "

# Row 4
$ws.Range("E4").Value = "  +0.04%  "

# Row 5
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "571.30"
$textForced.Add($c) | Out-Null
$ws.Range("E5").Value = "  -7.33%  "

# Row 6
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "185.80"
$textForced.Add($c) | Out-Null
$ws.Range("E6").Value = "  -4.33%  "

# Row 7
$ws.Range("D7").Value = "3.550.39"
$ws.Range("E7").Value = "  -4.42%  "

# Row 8
$ws.Range("E8").Value = "  -4.07%  "

# Row 9
$ws.Range("E9").Value = "  +0.19%  "

# Row 10
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "0.670"
$textForced.Add($c) | Out-Null
$ws.Range("E10").Value = "  -7.64%  "

# Row 11
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "0.149"
$textForced.Add($c) | Out-Null
$ws.Range("E11").Value = "  -7.49%  "

# Row 12
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "55.10"
$textForced.Add($c) | Out-Null
$ws.Range("E12").Value = "  -8.46%  "

# Row 13
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "0.0000262"
$textForced.Add($c) | Out-Null
$ws.Range("E13").Value = "  -9.68%  "

# Row 14
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "9.80"
$textForced.Add($c) | Out-Null
$ws.Range("E14").Value = "  -6.46%  "

# Row 15
$ws.Range("D15").Value = "4.120.58"
$ws.Range("E15").Value = "  -4.31%  "

# Row 16
$ws.Range("D16").Value = "3.551.64"
$ws.Range("E16").Value = "  -4.51%  "

# Row 17
$ws.Range("E17").Value = "  -1.74%  "

# Row 18
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "18.23"
$textForced.Add($c) | Out-Null
$ws.Range("E18").Value = "  -6.30%  "

# Row 19
$ws.Range("D19").Value = "66.405.01"
$ws.Range("E19").Value = "  -3.51%  "

# Row 20
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "12.05"
$textForced.Add($c) | Out-Null
$ws.Range("E20").Value = "  -6.97%  "

# Row 21
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "1.06"
$textForced.Add($c) | Out-Null
$ws.Range("E21").Value = "  -8.54%  "

# Row 22
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "388.19"
$textForced.Add($c) | Out-Null
$ws.Range("E22").Value = "  -5.84%  "

# Row 23
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "4.17"
$textForced.Add($c) | Out-Null
$ws.Range("E23").Value = "  -9.58%  "

# Row 24
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "85.22"
$textForced.Add($c) | Out-Null
$ws.Range("E24").Value = "  -5.66%  "

# Row 25
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "11.19"
$textForced.Add($c) | Out-Null
$ws.Range("E25").Value = "  -2.58%  "

# Row 26
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "2.91"
$textForced.Add($c) | Out-Null
$ws.Range("E26").Value = "  -6.39%  "

# Row 27
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "12.36"
$textForced.Add($c) | Out-Null
$ws.Range("E27").Value = "  -6.48%  "

# Row 28
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "6.04"
$textForced.Add($c) | Out-Null
$ws.Range("E28").Value = "  -0.27%  "

# Row 29
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "3.56"
$textForced.Add($c) | Out-Null
$ws.Range("E29").Value = "  -6.80%  "

# Row 30
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "8.82"
$textForced.Add($c) | Out-Null
$ws.Range("E30").Value = "  -9.23%  "

# Row 31
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "7.55"
$textForced.Add($c) | Out-Null
$ws.Range("E31").Value = "  -2.44%  "

# Row 32
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "30.78"
$textForced.Add($c) | Out-Null
$ws.Range("E32").Value = "  -6.40%  "

# Row 33
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "630.41"
$textForced.Add($c) | Out-Null
$ws.Range("E33").Value = "  -1.32%  "

# Row 34
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "12.12"
$textForced.Add($c) | Out-Null
$ws.Range("E34").Value = "  -5.07%  "

# Row 35
$ws.Range("E35").Value = "  -8.46%  "

# Row 36
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "63.28"
$textForced.Add($c) | Out-Null
$ws.Range("E36").Value = "  -6.20%  "

# Row 37
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "41.73"
$textForced.Add($c) | Out-Null
$ws.Range("E37").Value = "  -11.44%  "

# Row 38
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "0.401"
$textForced.Add($c) | Out-Null
$ws.Range("E38").Value = "  -3.51%  "

# Row 39
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "1.00"
$textForced.Add($c) | Out-Null
$ws.Range("E39").Value = "  +0.14%  "

# Row 40
$ws.Range("D40").Value = "0.0₃0745"
$ws.Range("E40").Value = "  -10.15%  "

# Row 41
$ws.Range("D41").Value = "3.133.18"
$ws.Range("E41").Value = "  +6.84%  "

# Row 42
$ws.Range("E42").Value = "  -6.21%  "

# Row 43
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "0.999"
$textForced.Add($c) | Out-Null
$ws.Range("E43").Value = "  -0.07%  "

# Row 46
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "0.0410"
$textForced.Add($c) | Out-Null
$ws.Range("E46").Value = "  -8.77%  "

# Row 47
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "0.130"
$textForced.Add($c) | Out-Null
$ws.Range("E47").Value = "  -7.07%  "

# Row 48
$ws.Range("E48").Value = "  -0.44%  "

# Row 49
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "139.24"
$textForced.Add($c) | Out-Null
$ws.Range("E49").Value = "  -4.44%  "

# Row 50
$ws.Range("E50").Value = "  -10.85%  "

# Row 51
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "2.73"
$textForced.Add($c) | Out-Null
$ws.Range("E51").Value = "  -1.63%  "

# Rows 44-45: Fetch.AI and ThetaToken swap ranking positions
$ws.Range("B44").Value = "Fetch.AI"
$ws.Range("C44").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "2.65"
$textForced.Add($c) | Out-Null
$ws.Range("E44").Value = "  +0.78%  "

$ws.Range("B45").Value = "ThetaToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "2.92"
$textForced.Add($c) | Out-Null
$ws.Range("E45").Value = "  -4.50%  "

# Restore default (General) formatting on every cell that was
# temporarily forced to text, so no stray style survives the edit.
foreach ($c in $textForced) {
    $c.ClearFormats()
}
